$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 3-11 have their data cyclically rotated: the content that used to be
# in row map[r] now lives in row r (only the columns below actually differ
# between rows 3-11; every other column already holds the same value in
# all nine rows, so leaving them untouched is equivalent to "rotating"
# them too - and avoids Excel's automatic text -> date coercion on the
# Y/Z/AA/AB "2022-11-14" / "00:00" text columns).
#   3<-5, 4<-6, 5<-7, 6<-8, 7<-9, 8<-3, 9<-10, 10<-11, 11<-4
$map = @{ 3 = 5; 4 = 6; 5 = 7; 6 = 8; 7 = 9; 8 = 3; 9 = 10; 10 = 11; 11 = 4 }

$cols = @("A", "B", "D", "E", "F", "G", "H", "L", "Q", "R", "AC")

# Snapshot the original values for every touched column/row before writing
# anything back (so later writes don't clobber values still to be read).
$orig = @{}
foreach ($col in $cols) {
    for ($r = 3; $r -le 11; $r++) {
        $orig["$col$r"] = $ws.Range("$col$r").Value()
    }
}

foreach ($col in $cols) {
    for ($destRow = 3; $destRow -le 11; $destRow++) {
        $srcRow = $map[$destRow]
        $ws.Range("$col$destRow").Value = $orig["$col$srcRow"]
    }
}
